$d = $word.ActiveDocument

# "added team information to softball details": the bullet describing the
# softball-league website gains a mention of team information alongside
# schedules and standings.
$old = "organize and display schedules and standings for the client’s"
$new = "organize and display schedules, standings, and team information for the client’s"

$found = $d.Content.Find.Execute(
    $old,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $new,
    2
)

if (-not $found) {
    throw "Could not find target text to replace."
}
